# Kaman new UI - header & footer
# The SCROLL_DOWN keyword step on the "TC14_Verify_AddToCart_from_Sear" sheet
# is replaced by a new TINY_SCROLL_DOWN keyword (cell B8), and that cell
# becomes the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B8").Value = "TINY_SCROLL_DOWN"
$ws.Range("B8").Select()
